$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P8: drop trailing "**" bullet from Foundation-of-Mathematics references ---
$ws.Range("P8").Value = '[''موحد، ضیاء، درآمدی به منطق جدید، انتشارات و آموزش انقلاب اسلامی ، 1379\n'', ''ناصر، بردجربان. مبانی و مقدمات علم ریاضی، مرکز نشر پرفسور حسابی، 1377\n'', ''سپهری ، محمد رضا. آشنائی با مبانی علوم ریاضی، مرکز نشر آینه، 1390\n'', ''Halmos, Paul R. I set theory. Courier Dover Publications, 2017.\n'', ''Stewart, Ian, and David Tall. The foundations of mathematics. OUP Oxford, 2015.\n'', ''Lin, Shwu-Yeng T., and You-Feng Lin. Set theory with applications. Mancorp Pub, 1985.\n'']'

# --- Row 9 <-> Row 10 swap: "Statistical Software" <-> "Probability and Statistics I" ---
# titles (fa/en/en-filename)
$ws.Range("B9").Value = 'آمار و احتمال ۱'
$ws.Range("C9").Value = 'Probability and Statistics I'
$ws.Range("D9").Value = 'Probability-and-Statistics-I'
$ws.Range("B10").Value = 'نرم‌افزارهای آماری و تحلیل داده‌ها'
$ws.Range("C10").Value = 'Statistical Software & data analysis'
$ws.Range("D10").Value = 'Statistical-Software-and-data-analysis'

# category: Probability moves from "required" to "base"
$ws.Range("E9").Value = 'پایه'

# unit count / unit type / hours - force text so numeric-looking strings stay as text
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = '3'
$ws.Range("I9").Style = "Normal"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = 'نظری'
$ws.Range("J9").Style = "Normal"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = '48'
$ws.Range("K9").Style = "Normal"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '1'
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = 'عملی'
$ws.Range("J10").Style = "Normal"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = '32'
$ws.Range("K10").Style = "Normal"

# exam column: Probability now has an exam (L9); Statistical Software has none (L10 cleared)
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = ' دارد'
$ws.Range("L9").Style = "Normal"
$ws.Range("L10").ClearContents()

# sar_fasl (course syllabus) and c_refs (references) swap between the two rows
$ws.Range("O9").Value = '[''احتمال: فضای احتمال، قوانین احتمال، مدل های احتمال، قوانین شمارش، احتمال شرطی، استقلال و قانون بیز. متغیرهای تصادفی: تعریف متغیر تصادفی، متغیرهای تصادفی (گسسته، پیوسته و آمیخته)، تابع توزیع، تابع جرم احتمال و تابع چگالی. امید ریاضی و گشتاورها: امید ریاضی، امید ریاضی تابعی از یک متغیر تصادفی. خواص و کاربردهای امید ریاضی، میانه و مد یک توزیع، واریانس و معیارهای پراکندگی دیگر، تقارن و چولگی، گشتاورهای یک متغیر تصادفی. تبدیل یک متغیر تصادفی: تابع مولد احتمال، تابع مولد گشتاور، ویژگیها و کاربردها. توزیع های استاندارد گسسته و پیوسته: برنولی، دوجمله ای، هندسی. فوق هندسی، دوجمله ای منفی، پواسون، یکنواخت (گسسته و پیوسته)، نمایی، گاما، بتا و نرمال. توزیع تابعی از یک متغیر تصادفی، تولید اعداد تصادفی از یک متغیر تصادفی، نامساوی های احتمالی: مارکوف، چپیشف، کشی و شوارتز.\n'']'
$ws.Range("O10").Value = '[''تعریف علم آمار، بیان کاربردهای آمار در سایر علوم، ارتباط و جایگاه علم آمار در علم داده\u200cها، جامعه آماری، انواع متغیرها، متغیرهای کمی  (پیوسته و گسسته) و متغیرهای رسته ای (اسمی و ترتیبی)، مقیاس ها، انواع داده\u200cها، داده\u200cهای یک متغیره و چند متغیره، داده\u200cهای سری زمانی، روشهای جمع آوری داده\u200cها، منابع داده\u200cهای ثبتی، معرفی مه داده\u200cها (کلان داده\u200cها) و منابع آنها، خلاصه کردن داده\u200cهای کمی یک متغیره با شاخص\u200cهای عددی،  انواع معیارهای مرکزیت و پراکندگی، چندکها، خلاصه کردن داده\u200cهای رسته\u200cای یک متغیره با شاخص\u200cهای عددی از قبیل نسبت\u200cها، انواع نمودارهای دیداری\u200cسازی داده\u200cهای یک متغیره کمی و رسته\u200cای، بافت نگار، منحنی فراوانی، توزیع تجربی، منحنی فراوانی نرمال،  نمودار ساقه و برگ، نمودار ستونی، نمودار جعبه ای، نمودار  سری زمانی و نمایش روند، تهیه جداول دو بعدی و چند بعدی برای خلاصه سازی داده\u200cهای چند متغیره رسته ای، نمودار دایره ای، انواع نمودارهای ستونی برای نمایش و توصیف داده\u200cهای چندمتغیره رسته ای، نمودار پراکنش برای sنمایش همبستگی در داده\u200cهای دو متغیره کمی و تفسیر آن، مقایسه نموداری دو یا چند سری زمانی، مقایسه میانگین دو یا چند جامعه و نمایش اختلافها با استفاده نمودار جعبه ای و سایر نمودارهای مناسب، اشاره ای به نمودار جریان، نمودارهای شبکه ای، درختی، نمودارهای گرمایی، قطبی و مارپیچی،  نقشه ای و جغرافیایی.\n'', ''**راهبردهای تدریس و یادگیری متناسب با محتوا و هدف:**\n'', ''برای پرورش  مهارت تحلیل داده\u200cها  و تسلط دانشجویان به نرم افزارهای آماری،  لازم است مدرس فعالیت های تعیین شده ای را در آزمایشگاه نرم افزار ارائه دهد و دانشجویان نیز فعالیت های تعریف شده ای را انجام دهند.  به طوری که نیمی از ساعات درس به آموزش تحلیل مقدماتی داده\u200cها  با استفاده از نرم افزار اکسل  اختصاص یابد.  دانشجویان نیز از طریق انجام پروژه های کلاسی با استفاده از داده\u200cهای واقعی و انجام تحلیل های مقدماتی گرافیکی و اکتشافی، در درس مشارکت داده شوند.  \n'', ''**ملزومات، تجهیزات و امکانات مورد نیاز برای ارائه:**  کلاس درس و آزمایشگاه رایانه مجهز\n'']'
$ws.Range("P9").Value = '[''بهبودیان، جواد. (1394). آمار واحتمال مقدماتی، انتشارات دانشگاه امام رضا، \n'', ''قهرمانی، سعید. (1397).  مبانی احتمال،  ترجمه شاهکار، غ. و بزرگنیا، ا.، ویرایش سوم، انتشارات دانشگاه صنعتی شریف.  \n'', ''راس، شلدون. (1389).  مبانی احتمال، ترجمه پارسیان، ا.  همدانی، ع.، ویرایش هشتم، انتشارات شیخ بهایی، اصفهان. \n'', ''هاگ، ر.و .تنیس، ل .آ .احتمال و استنباط آماری، ترجمه نوروز ایزددوستدار و حمید پزشک، انتشارات دانشگاه تهران،1۳۹۴\n'', ''Hogg, R. V. Tanis, E. and Zimmerman, D. Probability and Statistical Inference, 9th Edition, \n'', ''Pearson, 2013\n'', ''José Unpingco, Python for Probability, Statistics, and Machine Learning, Second Ed., 2019, Springer, Source Codes of the book: [Github Page](https://github.com/unpingco/Python-for-Probability-Statistics-and-Machine-Learning-2E)\n'', ''Grimmett, G. and Welsh D. Probability: an Introduction, 2nd Edition, OUP, 2014.\n'', ''Ghahramani, S., Fundamentals of Probability: with Stochastic Process, 3rd Edition, CRC Press,2014.\n'', ''**\n'']'
$ws.Range("P10").Value = '[''1\\. گوس، جی و منیروپ، د. (1401). آمار با نرم افزار JMP: نمودارها، آمار توصیفی و احتمال، ترجمه: بزرگنیا، ا.، و همکاران، انتشارات فردوسی مشهد. \n'', ''2\\. وناکت، ت. چ. و ووناکت، ر. ج. (1393). آمار مقدماتی، جلد 1 و2، ترجمه محمدرضا مشکانی، مرکز نشر دانشگاهی.\n'', ''2\\. Goos, P. and Meintrup, D. (2015).Statistics with JMP: graphs, descriptive statistics and probability. Wiley.\n'', ''3\\. Rahlf, T. (2019).Data Visualisation with R: 111 Eexamples. Springer Nature. \n'', ''4\\. Gerbing, D.W.(2020) R Visualisation: Derive meening from Data. Chpman and Hall/CRC.\n'', ''**\n'']'

# prerequisite references to "Probability and Statistics I" (now index/identity unaffected,
# but keep prereq text cells pointed at the literal course-name string)
$ws.Range("G28").Value = 'آمار و احتمال ۱'
$ws.Range("G34").Value = 'آمار و احتمال ۱'
$ws.Range("G35").Value = 'آمار و احتمال ۱'
$ws.Range("G36").Value = 'آمار و احتمال ۱'
$ws.Range("G43").Value = 'آمار و احتمال ۱'
$ws.Range("G64").Value = 'آمار و احتمال ۱'
